$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Angkatan"
$ws.Range("B2").Value = 20
$ws.Range("B3").Value = 21

$ws.Range("B4").Select()
